$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (6 rows x 2 columns), replacing the previous 32-row data
$data = @(
    @("https://www.google.co.in/maps/place/Geonla+(Brahamkhal)/data=!4m7!3m6!1s0x3908ea9cbae484bf:0x75d05b7637708a8!8m2!3d30.6992748!4d78.2984155!16s%2Fg%2F11b7grqfty!19sChIJv4TkupzqCDkRqAh3Y7cFXQc?authuser=0&hl=en&rclk=1", "Aboriginal art gallery"),
    @("https://www.google.co.in/maps/place/Kotdhar/data=!4m7!3m6!1s0x3908c14b22ca143f:0x375c00b346cac9cc!8m2!3d30.6749989!4d78.2300288!16s%2Fg%2F11s462p7_y!19sChIJPxTKIkvBCDkRzMnKRrMAXDc?authuser=0&hl=en&rclk=1", "Aboriginal art gallery"),
    @("https://www.google.co.in/maps/place/Banchaura/data=!4m7!3m6!1s0x3908c145b0381f9f:0x3f3c5dc0bbb9a3d5!8m2!3d30.6481765!4d78.215704!16s%2Fg%2F1hm3qgb9_!19sChIJnx84sEXBCDkR1aO5u8BdPD8?authuser=0&hl=en&rclk=1", "Aboriginal art gallery"),
    @("https://www.google.co.in/maps/place/Geonla+(Brahamkhal)/data=!4m7!3m6!1s0x3908ea9cbae484bf:0x75d05b7637708a8!8m2!3d30.6992748!4d78.2984155!16s%2Fg%2F11b7grqfty!19sChIJv4TkupzqCDkRqAh3Y7cFXQc?authuser=0&hl=en&rclk=1", "Aboriginal and Torres Strait Islander organization"),
    @("https://www.google.co.in/maps/place/Post+Office+Dunda/data=!4m7!3m6!1s0x3908eb1a2dfbdad1:0x599bb711a285e875!8m2!3d30.7069489!4d78.3475752!16s%2Fg%2F11cs383h0n!19sChIJ0dr7LRrrCDkRdeiFohG3m1k?authuser=0&hl=en&rclk=1", "Aboriginal and Torres Strait Islander organization"),
    @("https://www.google.co.in/maps/place/Banchaura/data=!4m7!3m6!1s0x3908c145b0381f9f:0x3f3c5dc0bbb9a3d5!8m2!3d30.6481765!4d78.215704!16s%2Fg%2F1hm3qgb9_!19sChIJnx84sEXBCDkR1aO5u8BdPD8?authuser=0&hl=en&rclk=1", "Aboriginal and Torres Strait Islander organization")
)

# Clear out the old used range first (previously rows 1-32, columns A-B)
$ws.UsedRange.Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
